$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-29 Monday", "2025-12-30 Tuesday"),
    @("23×13=299", "48×21=1008"),
    @("92×92=8464", "74×96=7104"),
    @("49×15=735", "58×87=5046"),
    @("36×30=1080", "80×87=6960"),
    @("46×47=2162", "13×65=845"),
    @("81×42=3402", "33×33=1089"),
    @("91×86=7826", "11×62=682"),
    @("39×48=1872", "35×20=700"),
    @("16×15=240", "62×70=4340"),
    @("92×99=9108", "46×78=3588"),
    @("32×82=2624", "84×85=7140"),
    @("54×50=2700", "69×57=3933"),
    @("37×80=2960", "90×49=4410"),
    @("67×20=1340", "91×69=6279"),
    @("58×88=5104", "69×73=5037"),
    @("24×86=2064", "79×22=1738"),
    @("14×67=938", "56×32=1792"),
    @("80×12=960", "26×62=1612"),
    @("91×31=2821", "52×21=1092"),
    @("65×82=5330", "49×84=4116"),
    @("32×53=1696", "51×37=1887"),
    @("20×33=660", "75×81=6075"),
    @("27×88=2376", "97×22=2134"),
    @("72×53=3816", "20×41=820"),
    @("65×90=5850", "88×60=5280")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
